$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 9478
$ws.Range("C2").Value = 9459
$ws.Range("D2").Value = 8396
$ws.Range("E2").Value = 0.887620255840998
$ws.Range("F2").Value = 0.8858408947035239
$ws.Range("G2").Value = 0.09668375143367978
$ws.Range("H2").Value = 0.08564642087330401
$ws.Range("I2").Value = 41289012.09786491
$ws.Range("J2").Value = 14448077.45109245
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 14448077.45109245
$ws.Range("M2").Value = 55737089.54895736
$ws.Range("N2").Value = 801867174.6472001
$ws.Range("O2").Value = 784167367.6432
$ws.Range("P2").Value = 0.01801804327187881
$ws.Range("Q2").Value = 0.01842473692129765
$ws.Range("B3").Value = 9666
$ws.Range("C3").Value = 9648
$ws.Range("D3").Value = 8564
$ws.Range("E3").Value = 0.8876451077943616
$ws.Range("F3").Value = 0.8859921373887855
$ws.Range("G3").Value = 0.1019731364521253
$ws.Range("H3").Value = 0.09034739712145676
$ws.Range("I3").Value = 48113730.21848053
$ws.Range("J3").Value = 17590471.60814859
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 17590471.60814859
$ws.Range("M3").Value = 65704201.82662913
$ws.Range("N3").Value = 838090574.9417281
$ws.Range("O3").Value = 820610398.917658
$ws.Range("P3").Value = 0.02098874767726823
$ws.Range("Q3").Value = 0.02143583804366786
$ws.Range("B4").Value = 9858
$ws.Range("C4").Value = 9838
$ws.Range("D4").Value = 8719
$ws.Range("E4").Value = 0.8862573693840211
$ws.Range("F4").Value = 0.8844593223777643
$ws.Range("G4").Value = 0.1062951315653054
$ws.Range("H4").Value = 0.09401372003630526
$ws.Range("I4").Value = 54572877.58804671
$ws.Range("J4").Value = 20500073.82904293
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 20500073.82904293
$ws.Range("M4").Value = 75072951.41708964
$ws.Range("N4").Value = 875624633.6923679
$ws.Range("O4").Value = 858175685.6864141
$ws.Range("P4").Value = 0.02341194278945468
$ws.Range("Q4").Value = 0.02388796859543496
$ws.Range("B5").Value = 10054
$ws.Range("C5").Value = 10027
$ws.Range("D5").Value = 8900
$ws.Range("E5").Value = 0.8876034706293009
$ws.Range("F5").Value = 0.8852198130097474
$ws.Range("G5").Value = 0.1095103629387546
$ws.Range("H5").Value = 0.0969407430032739
$ws.Range("I5").Value = 60694383.83057234
$ws.Range("J5").Value = 23227702.42889509
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 23227702.42889509
$ws.Range("M5").Value = 83922086.25946742
$ws.Range("N5").Value = 914123490.6623945
$ws.Range("O5").Value = 896637385.1994213
$ws.Range("P5").Value = 0.0254098080469016
$ws.Range("Q5").Value = 0.02590534681277986
$ws.Range("B6").Value = 10254
$ws.Range("C6").Value = 10233
$ws.Range("D6").Value = 9083
$ws.Range("E6").Value = 0.8876184892016027
$ws.Range("F6").Value = 0.8858006631558416
$ws.Range("G6").Value = 0.1088021642591696
$ws.Range("H6").Value = 0.09637702925356322
$ws.Range("I6").Value = 63966942.7227219
$ws.Range("J6").Value = 24492884.87752456
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 24492884.87752456
$ws.Range("M6").Value = 88459827.60024647
$ws.Range("N6").Value = 957320807.0430477
$ws.Range("O6").Value = 939728981.1698662
$ws.Range("P6").Value = 0.02558482454087431
$ws.Range("Q6").Value = 0.02606377516103998